$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 2.2
$ws.Range("I5").Value = 2.55

$ws.Range("G6").Value = 1.71

$ws.Range("G7").Value = 1.86

$ws.Range("G8").Value = 2.45
$ws.Range("I8").Value = 2.85
$ws.Range("N8").Value = 2.18
$ws.Range("P8").Value = 1.47
$ws.Range("Q8").Value = 2.35
$ws.Range("U8").Value = 11.25
$ws.Range("V8").Value = 9.75
$ws.Range("W8").Value = 26
$ws.Range("X8").Value = 23
$ws.Range("Y8").Value = 37
$ws.Range("Z8").Value = 7.3
$ws.Range("AA8").Value = 5.9
$ws.Range("AE8").Value = 7.4
$ws.Range("AF8").Value = 13.5
$ws.Range("AG8").Value = 10.75
$ws.Range("AH8").Value = 35
$ws.Range("AI8").Value = 28

$ws.Range("G11").Value = 1.67
$ws.Range("I11").Value = 5.25
$ws.Range("U11").Value = 7.5
$ws.Range("X11").Value = 13
$ws.Range("AG11").Value = 17

$ws.Range("N13").Value = 1.95
$ws.Range("O13").Value = 1.85

$ws.Range("G15").Value = 1.65
$ws.Range("H15").Value = 4
$ws.Range("I15").Value = 4.5
$ws.Range("AF15").Value = 23

$ws.Range("K21").Value = 10

$ws.Range("G26").Value = 1.85
$ws.Range("H26").Value = 3.6
$ws.Range("I26").Value = 4
$ws.Range("P26").Value = 1.33
$ws.Range("Q26").Value = 3.25
$ws.Range("R26").Value = 1.67
$ws.Range("S26").Value = 2.1
$ws.Range("T26").Value = 8.5
$ws.Range("U26").Value = 9.5
$ws.Range("W26").Value = 15
$ws.Range("Y26").Value = 23

$ws.Range("G27").Value = 1.95
$ws.Range("H27").Value = 3.75
$ws.Range("I27").Value = 3.4
$ws.Range("T27").Value = 10
$ws.Range("U27").Value = 11
$ws.Range("W27").Value = 19
$ws.Range("AA27").Value = 7.5
$ws.Range("AE27").Value = 15
$ws.Range("AF27").Value = 21
$ws.Range("AH27").Value = 41
